$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 11; this pushes the existing row 11 (and any
# rows below it) down by one, so the former row 11 becomes row 12. Excel's
# row insert carries the formatting (e.g. the date number format on column D)
# down from the row above, so the new row 11 already has the right styling.
$ws.Rows.Item(11).Insert()

# Populate the new row 11 with the new weekly price record.
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C11").Value = "Arica y Parinacota"
$ws.Range("D11").Value = 45132
$ws.Range("E11").Value = 15
$ws.Range("F11").Value = 100112030
$ws.Range("G11").Value = "Poroto granado"
$ws.Range("H11").Value = "Sin especificar"
$ws.Range("I11").Value = "Primera"
$ws.Range("J11").Value = 170
$ws.Range("K11").Value = 2200
$ws.Range("L11").Value = 2500
$ws.Range("M11").Value = 2359
$ws.Range("N11").Value = "$/kilo"
$ws.Range("O11").Value = "Región de Arica y Parinacota"
$ws.Range("P11").Value = 2359
$ws.Range("Q11").Value = 1
$ws.Range("R11").Value = "Hortaliza"
